# Append the latest statistics snapshot (2025-08-27 12:44 JST scrape run)
# as a new row at the bottom of the "統計" (statistics) sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("統計")

# Find the first empty row right after the existing data (row 3 -> row 4).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "2025-08-27T12:44:57.406498"
$ws.Cells.Item($newRow, 2).Value = 11
$ws.Cells.Item($newRow, 3).Value = "全案件リスト"
$ws.Cells.Item($newRow, 4).Value = 72.7
$ws.Cells.Item($newRow, 5).Value = 3
$ws.Cells.Item($newRow, 6).Value = 6
$ws.Cells.Item($newRow, 7).Value = 11
